$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.490.10'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '3.128.08'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.66'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.02'
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  -2.50%  '
$ws.Range("E8").Value = '  +3.64%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.125.79'
$ws.Range("E10").Value = '  +13.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.739'
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.203'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.79'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.56'
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").Value = '91.317.49'
$ws.Range("E16").Value = '  +1.01%  '
$ws.Range("E17").Value = '  +1.67%  '
$ws.Range("D18").Value = '3.132.04'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.74'
$ws.Range("E19").Value = '  -2.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.87'
$ws.Range("E20").Value = '  +4.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.84'
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '455.08'
$ws.Range("E22").Value = '  +2.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000201'
$ws.Range("E23").Value = '  -4.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.66'
$ws.Range("E25").Value = '  -3.19%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.37'
$ws.Range("E26").Value = '  -4.84%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.73'
$ws.Range("E27").Value = '  -2.26%  '
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.151'
$ws.Range("E28").Value = '  +41.73%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.317.67'
$ws.Range("E29").Value = '  +1.76%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.228'
$ws.Range("E31").Value = '  +3.38%  '
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.166'
$ws.Range("E32").Value = '  -4.89%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.34'
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.175'
$ws.Range("E34").Value = '  +10.93%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.23'
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.41'
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.99'
$ws.Range("E37").Value = '  +4.10%  '
$ws.Range("B38").Value = 'MantraDAO'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.92'
$ws.Range("E38").Value = '  -11.15%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '490.22'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.438'
$ws.Range("E41").Value = '  +5.46%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.37'
$ws.Range("E42").Value = '  -6.23%  '
$ws.Range("B43").Value = 'Binance-PegBSC-USD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.760'
$ws.Range("E43").Value = '  -23.97%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.704'
$ws.Range("E46").Value = '  +3.14%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '156.20'
$ws.Range("E47").Value = '  -1.90%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.92'
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("E50").Value = '  -3.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.06'
$ws.Range("E51").Value = '  -2.27%  '
